# Scoreboard.xlsx edit: add "Row1" (E) and "DP1" (F) score columns and a
# "Rep1" (G) total formula (=E+F) for every team row on the "Score" sheet,
# then update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-21 hold one team each. Fill in the new Row1 / DP1 raw-score
# columns with 0, matching the workbook's existing placeholder values.
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = 0   # column E - Row1
    $ws.Cells.Item($r, 6).Value = 0   # column F - DP1
}

# Column G ("Rep1") becomes the sum of the two new columns. Enter the
# first two rows individually (kept as standalone formulas) and then fill
# the remainder as one range so Excel groups them into a shared formula,
# mirroring how the workbook was actually authored.
$ws.Range("G2").Formula = "=E2+F2"
$ws.Range("G3").Formula = "=E3+F3"
$ws.Range("G4:G21").Formula = "=E4+F4"

# Move the active selection to L14, as recorded in the saved workbook.
$ws.Range("L14").Select() | Out-Null
